$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.78329501310442828
$ws.Range("D2").Value = 1.8466738095600985
$ws.Range("G2").Value = 0.10000405272976591
$ws.Range("H2").Value = 24.913461023337664
$ws.Range("I2").Value = 15.621224892810893
$ws.Range("J2").Value = 104.64310411888128
$ws.Range("K2").Value = 34.752449315596806
$ws.Range("L2").Value = 0.070536834606921867
$ws.Range("M2").Value = 15.596134821085466
$ws.Range("N2").Value = 11.194087049465692
$ws.Range("O2").Value = 76.375605484239955
$ws.Range("P2").Value = 26.696927526893106
$ws.Range("C3").Value = 0.98990221819709978
$ws.Range("D3").Value = 0.56283192302931706
$ws.Range("G3").Value = 0.1224401060204636
$ws.Range("H3").Value = 37.136149086996575
$ws.Range("I3").Value = 12.501858021838444
$ws.Range("J3").Value = 112.94276154022478
$ws.Range("K3").Value = 62.767671748572312
$ws.Range("L3").Value = 0.085856535583618054
$ws.Range("M3").Value = 23.056233999971766
$ws.Range("N3").Value = 8.906312758644706
$ws.Range("O3").Value = 81.976452938205483
$ws.Range("P3").Value = 48.074010965511434
$ws.Range("C4").Value = 0.96353449056561224
$ws.Range("D4").Value = 0.83754989478856323
$ws.Range("G4").Value = 0.10000031088873702
$ws.Range("H4").Value = 22.294682946419176
$ws.Range("I4").Value = 9.7950315256722238
$ws.Range("J4").Value = 109.05095718685831
$ws.Range("K4").Value = 35.557215284075696
$ws.Range("L4").Value = 0.072756222707026169
$ws.Range("M4").Value = 14.573465306989981
$ws.Range("N4").Value = 7.2379218238352507
$ws.Range("O4").Value = 81.93930147791346
$ws.Range("P4").Value = 27.782038659655161
$ws.Range("C5").Value = 0.9387437250173003
$ws.Range("D5").Value = 0.72308013428846418
$ws.Range("G5").Value = 0.10001263668525223
$ws.Range("H5").Value = 16.741674786246577
$ws.Range("I5").Value = 9.2406868972764755
$ws.Range("J5").Value = 71.585968678588344
$ws.Range("K5").Value = 23.550555029466061
$ws.Range("L5").Value = 0.071801392843886092
$ws.Range("M5").Value = 10.743010063252312
$ws.Range("N5").Value = 6.7391813753506025
$ws.Range("O5").Value = 53.124119847125691
$ws.Range("P5").Value = 18.263206024413417
$ws.Range("C6").Value = 0.9769215571354265
$ws.Range("D6").Value = 1.1171792241720118
$ws.Range("G6").Value = 0.10019816388634087
$ws.Range("H6").Value = 33.972219610744872
$ws.Range("I6").Value = 14.876994251487275
$ws.Range("J6").Value = 133.9527813567716
$ws.Range("K6").Value = 62.148068854647647
$ws.Range("L6").Value = 0.072370198585075957
$ws.Range("M6").Value = 21.983430639959064
$ws.Range("N6").Value = 10.91455344655621
$ws.Range("O6").Value = 99.968818316370971
$ws.Range("P6").Value = 48.357336232413111
$ws.Range("C7").Value = 0.99158087810887074
$ws.Range("D7").Value = 0.45670297870571192
$ws.Range("G7").Value = 0.10747930643896778
$ws.Range("H7").Value = 25.329025461355727
$ws.Range("I7").Value = 9.1195880696007201
$ws.Range("J7").Value = 109.76788687285288
$ws.Range("K7").Value = 52.75875054074892
$ws.Range("L7").Value = 0.077562246924497619
$ws.Range("M7").Value = 16.370811481433595
$ws.Range("N7").Value = 6.6849302841034275
$ws.Range("O7").Value = 81.853795691073046
$ws.Range("P7").Value = 41.031660537350803
$ws.Range("C8").Value = 0.96332099721190789
$ws.Range("D8").Value = 0.54798800107340739
$ws.Range("G8").Value = 0.10000883164966745
$ws.Range("H8").Value = 16.288500460874836
$ws.Range("I8").Value = 9.6071903359183342
$ws.Range("J8").Value = 81.456096172891918
$ws.Range("K8").Value = 30.840152337899877
$ws.Range("L8").Value = 0.069685477859098399
$ws.Range("M8").Value = 10.022765835410212
$ws.Range("N8").Value = 6.8008166690413026
$ws.Range("O8").Value = 58.768374786809055
$ws.Range("P8").Value = 23.546320064813138
$ws.Range("C9").Value = 0.98237650834681167
$ws.Range("D9").Value = 0.62011373959013427
$ws.Range("G9").Value = 0.2235730917971317
$ws.Range("H9").Value = 28.617916164324047
$ws.Range("I9").Value = 11.827221251124925
$ws.Range("J9").Value = 106.81211596992512
$ws.Range("K9").Value = 27.914289772664617
$ws.Range("L9").Value = 0.15077770811999974
$ws.Range("M9").Value = 16.801817438367742
$ws.Range("N9").Value = 8.0978546129301865
$ws.Range("O9").Value = 74.668988612271193
$ws.Range("P9").Value = 20.997657198006159
$ws.Range("C10").Value = 0.96798126390294659
$ws.Range("D10").Value = 0.58434587666726401
$ws.Range("G10").Value = 0.10001559407935941
$ws.Range("H10").Value = 18.085118965196308
$ws.Range("I10").Value = 8.7066131941273479
$ws.Range("J10").Value = 68.751400329459344
$ws.Range("K10").Value = 24.226209772259221
$ws.Range("L10").Value = 0.070867058338375935
$ws.Range("M10").Value = 11.394170778050862
$ws.Range("N10").Value = 6.2675207339667924
$ws.Range("O10").Value = 50.395601868316106
$ws.Range("P10").Value = 18.65487186801095
$ws.Range("C11").Value = 0.75097529748476766
$ws.Range("D11").Value = 0.40883560796923629
$ws.Range("G11").Value = 29.85426286542269
$ws.Range("H11").Value = 6.3415887275175846
$ws.Range("I11").Value = 4.736850224554324
$ws.Range("J11").Value = 52.034610543388993
$ws.Range("K11").Value = 7.046652307766065
$ws.Range("L11").Value = 20.728423224012886
$ws.Range("M11").Value = 3.8825061583100964
$ws.Range("N11").Value = 3.3411619793822029
$ws.Range("O11").Value = 37.414432587737224
$ws.Range("P11").Value = 5.3707655081906838
$ws.Range("C12").Value = 0.99683168511846221
$ws.Range("D12").Value = 0.20871156917643785
$ws.Range("G12").Value = 0.16135302220977976
$ws.Range("H12").Value = 17.931653379531369
$ws.Range("I12").Value = 7.3795891155854596
$ws.Range("J12").Value = 68.480224050215668
$ws.Range("K12").Value = 24.151418513943053
$ws.Range("L12").Value = 0.11221998558574356
$ws.Range("M12").Value = 11.004649430460649
$ws.Range("N12").Value = 5.2141056405976993
$ws.Range("O12").Value = 49.318752951040167
$ws.Range("P12").Value = 18.422653998899335
$ws.Range("C13").Value = 0.95943866745095074
$ws.Range("D13").Value = 1.0873018463492452
$ws.Range("G13").Value = 0.10006922382184548
$ws.Range("H13").Value = 27.423475529793642
$ws.Range("I13").Value = 9.5928219279264368
$ws.Range("J13").Value = 97.066218998087578
$ws.Range("K13").Value = 48.230761019648163
$ws.Range("L13").Value = 0.070399575860043459
$ws.Range("M13").Value = 17.104737924509298
$ws.Range("N13").Value = 6.8563291564505588
$ws.Range("O13").Value = 70.671516265318317
$ws.Range("P13").Value = 37.001538243916322
$ws.Range("C14").Value = 0.97326813102549958
$ws.Range("D14").Value = 0.39611002820300767
$ws.Range("G14").Value = 0.10001208672522731
$ws.Range("H14").Value = 13.362410821891572
$ws.Range("I14").Value = 5.883613364416461
$ws.Range("J14").Value = 67.517462252747904
$ws.Range("K14").Value = 28.363171934830547
$ws.Range("L14").Value = 0.068643736281863918
$ws.Range("M14").Value = 8.0471345448732627
$ws.Range("N14").Value = 4.1017525489001905
$ws.Range("O14").Value = 48.012379008799137
$ws.Range("P14").Value = 21.500334226107103
$ws.Range("C15").Value = 0.82714350169032713
$ws.Range("D15").Value = 0.96125685126048033
$ws.Range("G15").Value = 29.987713859243719
$ws.Range("H15").Value = 19.562202927333093
$ws.Range("I15").Value = 13.15852159410127
$ws.Range("J15").Value = 83.044575098429405
$ws.Range("K15").Value = 11.958649945618919
$ws.Range("L15").Value = 20.086201534082424
$ws.Range("M15").Value = 11.37096153706103
$ws.Range("N15").Value = 8.9456810187482905
$ws.Range("O15").Value = 57.665418072430604
$ws.Range("P15").Value = 8.9704385811837923
